$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.430.46"
$ws.Range("E2").Value = "  +1.71%  "

$ws.Range("D3").Value = "1.669.16"
$ws.Range("E3").Value = "  +2.05%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.33"
$ws.Range("E5").Value = "  +1.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3969"
$ws.Range("E7").Value = "  +2.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3927"
$ws.Range("E8").Value = "  +2.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.39"
$ws.Range("E9").Value = "  +6.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.394"
$ws.Range("E10").Value = "  +4.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08566"
$ws.Range("E12").Value = "  +2.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.53"
$ws.Range("E13").Value = "  +4.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.300"
$ws.Range("E14").Value = "  +4.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.945"
$ws.Range("E15").Value = "  +7.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001335"
$ws.Range("E16").Value = "  +5.41%  "

$ws.Range("D17").Value = "1.662.64"
$ws.Range("E17").Value = "  +0.73%  "

$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07033"
$ws.Range("E19").Value = "  +2.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.66"
$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.994"
$ws.Range("E21").Value = "  +2.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("E23").Value = "  +2.39%  "

$ws.Range("D24").Value = "24.437.67"
$ws.Range("E24").Value = "  +1.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.480"
$ws.Range("E25").Value = "  +6.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.066"
$ws.Range("E26").Value = "  +15.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.54"
$ws.Range("E27").Value = "  +1.58%  "

$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.454"
$ws.Range("E29").Value = "  +2.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "142.53"
$ws.Range("E30").Value = "  +2.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.033"
$ws.Range("E31").Value = "  -6.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.547"
$ws.Range("E32").Value = "  +5.60%  "

$ws.Range("D33").Value = "1.845.22"
$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("E34").Value = "  +13.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03105"
$ws.Range("E35").Value = "  +8.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08281"
$ws.Range("E36").Value = "  +4.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.928"
$ws.Range("E37").Value = "  +2.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.21"
$ws.Range("E38").Value = "  +14.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2766"
$ws.Range("E39").Value = "  +4.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09270"
$ws.Range("E40").Value = "  +1.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7716"
$ws.Range("E41").Value = "  +3.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.76"
$ws.Range("E42").Value = "  +6.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.444"
$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.53"
$ws.Range("E44").Value = "  +3.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7121"
$ws.Range("E45").Value = "  +4.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.548"
$ws.Range("E46").Value = "  +4.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.129"
$ws.Range("E47").Value = "  +1.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9995"
$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08446"
$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.94"
$ws.Range("E50").Value = "  +4.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.269"
$ws.Range("E51").Value = "  +2.25%  "
